$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24-78 down to 25-79.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new record.
$ws.Range("A24").Value = 5
$ws.Range("B24").Value = "Macroferia Regional de Talca"
$ws.Range("C24").Value = "Maule"
$ws.Range("D24").Value = 44565
$ws.Range("E24").Value = 7
$ws.Range("F24").Value = 100112001
$ws.Range("G24").Value = "Berenjena"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 150
$ws.Range("K24").Value = 9000
$ws.Range("L24").Value = 9000
$ws.Range("M24").Value = 9000
$ws.Range("N24").Value = "$/caja 50 unidades"
$ws.Range("O24").Value = "Región del Maule"
$ws.Range("P24").Value = 180
$ws.Range("Q24").Value = 50
$ws.Range("R24").Value = "Hortaliza"
